$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "text-number" formatting convention used throughout column D
# (prices such as "27.940.22" or "1.005" are stored as literal text, not numeric values),
# so force Text format before writing the new price strings to avoid Excel auto-converting
# them into real numbers (which would strip trailing zeros / change precision).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.914.35"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.748.15"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "333.92"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.3856"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "0.3374"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").Value = "1.111"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "6.148"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "1.748.86"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "7.078"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "0.00001055"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "0.06611"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "79.05"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "16.73"
$ws.Range("E21").Value = "  -3.39%  "
$ws.Range("D22").Value = "6.170"
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").Value = "27.923.84"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "11.62"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").Value = "2.403"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "153.60"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "19.76"
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("E28").Value = "  -5.09%  "
$ws.Range("D29").Value = "1.950.50"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").Value = "1.285"
$ws.Range("E30").Value = "  -9.90%  "
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").Value = "4.024"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").Value = "5.794"
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("D34").Value = "0.08780"
$ws.Range("D35").Value = "12.12"
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("D36").Value = "1.536"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").Value = "0.6509"
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "5.130"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02273"
$ws.Range("E39").Value = "  -5.62%  "
$ws.Range("D40").Value = "0.06108"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").Value = "0.2090"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").Value = "7.947"
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("D45").Value = "13.68"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").Value = "3.825"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "0.6012"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").Value = "126.79"
$ws.Range("E48").Value = "  -2.91%  "
$ws.Range("D49").Value = "1.988"
$ws.Range("E49").Value = "  -4.01%  "
$ws.Range("D50").Value = "1.167"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").Value = "1.106"
$ws.Range("E51").Value = "  +4.11%  "
